$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated betting odds / count values for rows 3, 4, 5, 6 and 8
# as scraped on 2024-10-03 (FlashScore weekly games refresh).
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 10
$ws.Range("O3").Value = 1.36
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 2.15
$ws.Range("R3").Value = 1.67
$ws.Range("G4").Value = 1.67
$ws.Range("H4").Value = 3.25
$ws.Range("I4").Value = 6.25
$ws.Range("J4").Value = 2.4
$ws.Range("K4").Value = 1.91
$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 5.5
$ws.Range("O4").Value = 1.57
$ws.Range("P4").Value = 2.25
$ws.Range("Z4").Value = 12
$ws.Range("AC4").Value = 5.5
$ws.Range("AF4").Value = 101
$ws.Range("AK4").Value = 51
$ws.Range("AN4").Value = 3.4
$ws.Range("AO4").Value = 9.5
$ws.Range("AW4").Value = 7
$ws.Range("AZ4").Value = 151
$ws.Range("O5").Value = 1.25
$ws.Range("P5").Value = 3.75
$ws.Range("Q5").Value = 1.85
$ws.Range("R5").Value = 1.95
$ws.Range("S5").Value = 1.36
$ws.Range("T5").Value = 3
$ws.Range("W5").Value = 7
$ws.Range("Y5").Value = 8.5
$ws.Range("AT5").Value = 3
$ws.Range("I6").Value = 1.53
$ws.Range("J6").Value = 4.75
$ws.Range("L6").Value = 2
$ws.Range("O6").Value = 1.13
$ws.Range("P6").Value = 5.5
$ws.Range("Q6").Value = 1.44
$ws.Range("R6").Value = 2.63
$ws.Range("X6").Value = 34
$ws.Range("AG6").Value = 10
$ws.Range("AI6").Value = 9
$ws.Range("AJ6").Value = 12
$ws.Range("AM6").Value = 126
$ws.Range("AX6").Value = 7.5
$ws.Range("G8").Value = 2.9
$ws.Range("I8").Value = 2.25
$ws.Range("J8").Value = 3.4
$ws.Range("L8").Value = 2.88
$ws.Range("S8").Value = 1.3
$ws.Range("T8").Value = 3.4
$ws.Range("W8").Value = 12
$ws.Range("X8").Value = 17
$ws.Range("Y8").Value = 11
$ws.Range("AH8").Value = 12
$ws.Range("AI8").Value = 9
$ws.Range("AJ8").Value = 21
$ws.Range("AL8").Value = 21
$ws.Range("AQ8").Value = 51
$ws.Range("AT8").Value = 3.4
$ws.Range("AW8").Value = 4.5
$ws.Range("AX8").Value = 12
$ws.Range("BB8").Value = 101
